# m1-8-10.xml: p. 2 of 3
#
# This edit:
#   1. Inserts two blank rows above row 13 on the "Transcriptions" sheet
#      (pushing the M45.. M113 block, plus the trailing note row, down by two).
#   2. Adds a link (column D) for the existing "London Bridge" annotation row,
#      and appends a brand-new annotation row for "James I of England and VI
#      of Scotland" on the "Annotations" sheet.
#   3. Leaves the "Annotations" sheet as the active/selected tab, with the
#      selection sitting on the newly-added row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Transcriptions sheet: insert two rows above row 13.
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transcriptions")

$insertRows = $wsTrans.Range("A13:H14")
$insertRows.EntireRow.Insert()

# Selection on the Transcriptions sheet, matching the saved view state.
$wsTrans.Range("B18").Select()

# ---------------------------------------------------------------------------
# 2. Annotations sheet: fill in the missing link for "London Bridge" and add
#    a new row for "James I of England and VI of Scotland".
# ---------------------------------------------------------------------------
$wsAnnot = $wb.Worksheets.Item("Annotations")

$wsAnnot.Range("D236").Value = "../resources/annotations.xml#pla-lobr"

$wsAnnot.Range("A237").Value = "James I of England and VI of Scotland"
$wsAnnot.Range("B237").Value = "Person"
$wsAnnot.Range("C237").Value = "psn-kjam"
$wsAnnot.Range("D237").Value = "../resources/annotations.xml#psn-kjam"

$wsAnnot.Range("D237").Select()

# The Annotations sheet ends up as the active tab.
$wsAnnot.Activate()
